$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New player / position / team data that replaces the existing A2:C19 table.
$data = @(
    ,@('Collin Sexton', 'PG,SG', 'Utah Jazz')
    ,@('Malik Monk', 'PG,SG,SF', 'Sacramento Kings')
    ,@('Tyus Jones', 'PG', 'Phoenix Suns')
    ,@('Duncan Robinson', 'SG,SF', 'Miami Heat')
    ,@('Devin Vassell', 'SG,SF', 'San Antonio Spurs')
    ,@('Tari Eason', 'SF,PF', 'Houston Rockets')
    ,@('Onyeka Okongwu', 'PF,C', 'Atlanta Hawks')
    ,@('Naz Reid', 'PF,C', 'Minnesota Timberwolves')
    ,@('Deandre Ayton', 'C', 'Portland Trail Blazers')
    ,@('Isaiah Hartenstein', 'C', 'Oklahoma City Thunder')
    ,@('Derrick White', 'PG,SG', 'Boston Celtics')
    ,@('Coby White', 'PG,SG', 'Chicago Bulls')
    ,@('Damian Lillard', 'PG', 'Milwaukee Bucks')
    ,@('Cade Cunningham', 'PG,SG', 'Detroit Pistons')
    ,@('Klay Thompson', 'SG,SF', 'Dallas Mavericks')
    ,@('Anthony Davis', 'PF,C', 'Los Angeles Lakers')
    ,@('LaMelo Ball', 'PG,SG', 'Charlotte Hornets')
    ,@('Julius Randle', 'PF,C', 'Minnesota Timberwolves')
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
